$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 45, pushing the old row 45 down to row 46.
$ws.Rows.Item(45).Insert()

# New row 45 gets the latest data point (week of 44448), copying the
# unchanged columns (A, B, C) from the row that was pushed down and
# updating the data columns with the new values.
$ws.Cells.Item(45, 1).Value = 10
$ws.Cells.Item(45, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(45, 3).Value = "La Araucanía"
$ws.Cells.Item(45, 4).Value = 44448
$ws.Cells.Item(45, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(45, 5).Value = 9
$ws.Cells.Item(45, 6).Value = 100112035
$ws.Cells.Item(45, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(45, 8).Value = "Sin especificar"
$ws.Cells.Item(45, 9).Value = "Primera"
$ws.Cells.Item(45, 10).Value = 15
$ws.Cells.Item(45, 11).Value = 25000
$ws.Cells.Item(45, 12).Value = 25000
$ws.Cells.Item(45, 13).Value = 25000
$ws.Cells.Item(45, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(45, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(45, 16).Value = 2500
$ws.Cells.Item(45, 17).Value = 10
$ws.Cells.Item(45, 18).Value = "Hortaliza"

$wb.Save()
